$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "取得日時" (retrieved at) timestamp in column A for rows 2-19
# from 2025-11-25 01:20:23 to 2025-11-25 01:50:54 (appended at 01:50 JST).
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-25 01:50:54"
}
